$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.613.47'
$ws.Range("E2").Value = '  +4.83%  '

$ws.Range("D3").Value = '2.268.95'

$ws.Range("E4").Value = '  +0.07%  '

$ws.Range("D5").Value = "'230.68"
$ws.Range("E5").Value = '  +0.18%  '

$ws.Range("D6").Value = "'0.626"
$ws.Range("E6").Value = '  +1.05%  '

$ws.Range("D7").Value = "'63.22"
$ws.Range("E7").Value = '  +5.86%  '

$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("D9").Value = "'0.427"
$ws.Range("E9").Value = '  +6.28%  '

$ws.Range("D10").Value = "'0.106"
$ws.Range("E10").Value = '  +18.81%  '

$ws.Range("D11").Value = "'57.24"
$ws.Range("E11").Value = '  -0.77%  '

$ws.Range("D12").Value = "'25.83"
$ws.Range("E12").Value = '  +15.30%  '

$ws.Range("D13").Value = "'0.103"
$ws.Range("E13").Value = '  +0.03%  '

$ws.Range("D14").Value = '2.609.21'
$ws.Range("E14").Value = '  +2.27%  '

$ws.Range("D15").Value = "'15.64"
$ws.Range("E15").Value = '  +1.37%  '

$ws.Range("D16").Value = "'5.89"
$ws.Range("E16").Value = '  +4.85%  '

$ws.Range("D17").Value = "'0.820"
$ws.Range("E17").Value = '  +3.00%  '

$ws.Range("D18").Value = '2.275.90'
$ws.Range("E18").Value = '  +1.69%  '

$ws.Range("D19").Value = '43.503.76'
$ws.Range("E19").Value = '  +4.59%  '

$ws.Range("E20").Value = '  +11.66%  '

$ws.Range("D21").Value = "'73.27"
$ws.Range("E21").Value = '  +1.50%  '

$ws.Range("E22").Value = '  +0.00%  '

$ws.Range("D23").Value = "'249.31"
$ws.Range("E23").Value = '  +1.74%  '

$ws.Range("E24").Value = '  +0.17%  '

$ws.Range("D25").Value = "'2.48"
$ws.Range("E25").Value = '  +5.13%  '

$ws.Range("E26").Value = '  +0.82%  '

$ws.Range("D27").Value = "'9.82"
$ws.Range("E27").Value = '  +1.11%  '

$ws.Range("D28").Value = "'171.57"
$ws.Range("E28").Value = '  +1.61%  '

$ws.Range("D29").Value = "'20.91"
$ws.Range("E29").Value = '  +5.70%  '

$ws.Range("D30").Value = "'0.137"
$ws.Range("E30").Value = '  -2.09%  '

$ws.Range("D31").Value = "'1.43"
$ws.Range("E31").Value = '  +1.93%  '

$ws.Range("E32").Value = '  +11.11%  '

$ws.Range("E33").Value = '  +0.86%  '

$ws.Range("D34").Value = "'0.0684"
$ws.Range("E34").Value = '  +5.15%  '

$ws.Range("D35").Value = "'5.04"
$ws.Range("E35").Value = '  +1.67%  '

$ws.Range("D36").Value = "'4.70"
$ws.Range("E36").Value = '  +0.92%  '

$ws.Range("E37").Value = '  +4.41%  '

$ws.Range("D38").Value = "'3.79"
$ws.Range("E38").Value = '  +5.97%  '

$ws.Range("D39").Value = "'2.32"
$ws.Range("E39").Value = '  -2.07%  '

$ws.Range("D40").Value = "'0.0247"
$ws.Range("E40").Value = '  +4.24%  '

$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = '  +0.10%  '

$ws.Range("D42").Value = "'8.37"
$ws.Range("E42").Value = '  -2.88%  '

$ws.Range("B43").Value = 'InjectiveProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D43").Value = "'17.22"
$ws.Range("E43").Value = '  +4.14%  '

$ws.Range("B44").Value = 'Cronos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D44").Value = "'0.0960"
$ws.Range("E44").Value = '  -0.31%  '

$ws.Range("B45").Value = 'FTXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D45").Value = "'4.42"
$ws.Range("E45").Value = '  +1.23%  '

$ws.Range("D46").Value = "'1.20"
$ws.Range("E46").Value = '  -0.71%  '

$ws.Range("D47").Value = "'97.22"
$ws.Range("E47").Value = '  -0.42%  '

$ws.Range("B48").Value = 'Celestia'
$ws.Range("C48").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D48").Value = "'10.26"
$ws.Range("E48").Value = '  +18.75%  '

$ws.Range("D49").Value = '1.474.50'
$ws.Range("E49").Value = '  +0.35%  '

$ws.Range("D50").Value = "'2.33"
$ws.Range("E50").Value = '  +4.29%  '

$ws.Range("B51").Value = 'ARBITRUM'
$ws.Range("C51").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D51").Value = "'1.07"
$ws.Range("E51").Value = '  -0.12%  '
